$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.965.98"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.018.34"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'225.75"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("D8").Value = "'54.85"
$ws.Range("E8").Value = "  -4.66%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  -3.95%  "
$ws.Range("D12").Value = "2.318.92"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "'14.27"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").Value = "'20.33"
$ws.Range("E14").Value = "  -2.76%  "
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").Value = "2.008.52"
$ws.Range("E17").Value = "  -2.61%  "
$ws.Range("D18").Value = "36.881.94"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'6.24"
$ws.Range("E19").Value = "  +5.17%  "
$ws.Range("D20").Value = "'68.84"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'225.84"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.41"
$ws.Range("E24").Value = "  +2.17%  "
$ws.Range("E25").Value = "  -7.63%  "
$ws.Range("D26").Value = "'165.55"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("E27").Value = "  -4.17%  "
$ws.Range("E28").Value = "  -5.30%  "
$ws.Range("D29").Value = "'18.69"
$ws.Range("E29").Value = "  -3.71%  "
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("D32").Value = "'4.50"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("E34").Value = "  -4.03%  "
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  -4.90%  "
$ws.Range("D39").Value = "'5.29"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "'17.14"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  -5.24%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.482.14"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "'95.36"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").Value = "'0.0927"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  -5.21%  "
$ws.Range("E46").Value = "  -6.28%  "
$ws.Range("D47").Value = "'7.35"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  -3.44%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").Value = "2.206.96"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("E51").Value = "  -9.38%  "
